# Auto-generated edit script: updates the cryptos price table
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.919.50"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.49"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.93%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7420"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.70"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3148"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07165"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.76"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08417"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.26%  "

$ws.Range("E12").Value = "  -2.29%  "

$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.872.46"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -7.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.59"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.73%  "

$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.112"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.52%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.910.88"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.60"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.48"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007819"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9994"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.124.97"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -6.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.985"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9993"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1561"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.309"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.79"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.63"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.041"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.483"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.46%  "

$ws.Range("E31").Value = "  +2.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.531"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.263"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05334"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.239"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.70%  "

$ws.Range("E36").Value = "  +0.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9970"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.698"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01953"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.752"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4498"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.114.29"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.064"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.29"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8585"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.40%  "

$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.16"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.674"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.080"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.844"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.022.87"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.57%  "

